# Add a new "Greece" test-data sheet, cloned from the existing "Croatia"
# sheet (same layout/styles), filled in with the Greece-specific values,
# and placed as the last (now active) tab - mirrors the commit
# "Test data for Greece Market".

$wb = $excel.ActiveWorkbook

# Template sheet to clone: the last country sheet ("Croatia").
$template = $wb.Worksheets.Item("Croatia")
$template.Activate()
$template.Cells.Select()

# Copy it to the end of the workbook (After:=$template, Before:=$null).
$template.Copy($null, $template)

# The copy becomes the active sheet/tab.
$greece = $wb.ActiveSheet
$greece.Name = "Greece"

# Fill in the Greece-specific values (order matters for shared-string
# insertion order: model/Wg ref first, then the market name).
$greece.Range("B4").Value = "NGC-4119/T3187/T3189"
$greece.Range("B2").Value = "Greece Market"

# Leave the cursor where the author left it on the new sheet.
$greece.Range("D14").Select()
